{"js": "// Replace each three-digit-by-one-digit multiplication expression in the\n// practice table with its new value, per the commit's regenerated numbers.\n// Old -> New pairs (each old value is unique in the document).\nconst replacements = [\n  [\"818\u00d75=\", \"643\u00d76=\"],\n  [\"773\u00d75=\", \"425\u00d77=\"],\n  [\"354\u00d76=\", \"116\u00d72=\"],\n  [\"304\u00d79=\", \"779\u00d72=\"],\n  [\"531\u00d73=\", \"937\u00d77=\"],\n  [\"855\u00d73=\", \"717\u00d74=\"],\n  [\"577\u00d78=\", \"570\u00d72=\"],\n  [\"585\u00d72=\", \"417\u00d76=\"],\n  [\"976\u00d74=\", \"136\u00d79=\"],\n  [\"313\u00d78=\", \"360\u00d78=\"],\n  [\"908\u00d72=\", \"659\u00d74=\"],\n  [\"995\u00d78=\", \"930\u00d72=\"],\n  [\"594\u00d75=\", \"199\u00d77=\"],\n  [\"975\u00d72=\", \"545\u00d76=\"],\n  [\"796\u00d73=\", \"658\u00d79=\"],\n  [\"280\u00d77=\", \"929\u00d74=\"],\n  [\"134\u00d73=\", \"435\u00d78=\"],\n  [\"794\u00d73=\", \"214\u00d78=\"],\n  [\"538\u00d77=\", \"867\u00d76=\"],\n  [\"355\u00d79=\", \"339\u00d79=\"],\n  [\"728\u00d76=\", \"428\u00d76=\"],\n  [\"980\u00d78=\", \"518\u00d77=\"],\n  [\"588\u00d79=\", \"908\u00d72=\"],\n  [\"964\u00d77=\", \"846\u00d76=\"],\n  [\"219\u00d73=\", \"458\u00d79=\"],\n];\n\nconst body = context.document.body;\n\n// Find every occurrence of each OLD string first (against the document's\n// original content), then apply the text replacements. Doing the lookups\n// up front avoids any chance that a newly written value (e.g. \"908\u00d72=\",\n// which is both an old value earlier in the table and a new value later)\n// gets matched again by a later search.\nconst searchResults = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\n\nsearchResults.forEach((result) => {\n  result.load(\"items\");\n});\n\nawait context.sync();\n\nsearchResults.forEach((result, i) => {\n  const [, newText] = replacements[i];\n  if (result.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${replacements[i][0]}\"`);\n  }\n  result.items.forEach((range) => {\n    range.insertText(newText, \"Replace\");\n  });\n});\n\nawait context.sync();\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression in the\n# practice table with its new value, per the commit's regenerated numbers.\n# Old -> New pairs (each old value is unique in the document). The pairs are\n# listed in document order, which matters because one new value (\"908\u00d72=\")\n# happens to equal another pair's old value earlier in the list - running the\n# Find/Replace calls in this order guarantees the earlier cell is already\n# rewritten before the later cell reintroduces that text.\n$d = $word.ActiveDocument\n\n$replacements = @()\n$replacements += ,@(\"818\u00d75=\", \"643\u00d76=\")\n$replacements += ,@(\"773\u00d75=\", \"425\u00d77=\")\n$replacements += ,@(\"354\u00d76=\", \"116\u00d72=\")\n$replacements += ,@(\"304\u00d79=\", \"779\u00d72=\")\n$replacements += ,@(\"531\u00d73=\", \"937\u00d77=\")\n$replacements += ,@(\"855\u00d73=\", \"717\u00d74=\")\n$replacements += ,@(\"577\u00d78=\", \"570\u00d72=\")\n$replacements += ,@(\"585\u00d72=\", \"417\u00d76=\")\n$replacements += ,@(\"976\u00d74=\", \"136\u00d79=\")\n$replacements += ,@(\"313\u00d78=\", \"360\u00d78=\")\n$replacements += ,@(\"908\u00d72=\", \"659\u00d74=\")\n$replacements += ,@(\"995\u00d78=\", \"930\u00d72=\")\n$replacements += ,@(\"594\u00d75=\", \"199\u00d77=\")\n$replacements += ,@(\"975\u00d72=\", \"545\u00d76=\")\n$replacements += ,@(\"796\u00d73=\", \"658\u00d79=\")\n$replacements += ,@(\"280\u00d77=\", \"929\u00d74=\")\n$replacements += ,@(\"134\u00d73=\", \"435\u00d78=\")\n$replacements += ,@(\"794\u00d73=\", \"214\u00d78=\")\n$replacements += ,@(\"538\u00d77=\", \"867\u00d76=\")\n$replacements += ,@(\"355\u00d79=\", \"339\u00d79=\")\n$replacements += ,@(\"728\u00d76=\", \"428\u00d76=\")\n$replacements += ,@(\"980\u00d78=\", \"518\u00d77=\")\n$replacements += ,@(\"588\u00d79=\", \"908\u00d72=\")\n$replacements += ,@(\"964\u00d77=\", \"846\u00d76=\")\n$replacements += ,@(\"219\u00d73=\", \"458\u00d79=\")\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $oldText,  # FindText\n        $false,    # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $newText,  # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Could not find text to replace: '$oldText'\"\n    }\n}\n\nWrite-Output \"done\"\n"}
